$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 3")
$tr = $sh.TextFrame.TextRange

$lines = @(
  "Docker",
  "`tImages",
  "`tContainers",
  "`tBuild a Docker image of your own.",
  "`tCreate a docker image of our application.",
  "`tPush it to Docker Hub",
  "`tLaunch it on an AWS EC2 ubuntu docker instance.",
  "`tCI with the Company/Employee microservices that we built.",
  "",
  "Troubleshooting Spring Boot application"
)
$tr.Text = [string]::Join("`r", $lines)

# Split "built." into its own run (still bold) so it can carry its own
# run properties, matching the authored edit.
$full = $tr.Text
$anchor = "`tCI with the Company/Employee microservices that we built."
$anchorIdx = $full.IndexOf($anchor)
$builtIdx = $anchorIdx + $anchor.Length - "built.".Length
$builtRun = $tr.Characters($builtIdx + 1, 6)
$builtRun.Font.Bold = $false
$builtRun.Font.Bold = $true
